$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.697.65'
$ws.Range("E2").Value = '  +2.72%  '
$ws.Range("D3").Value = '3.479.84'
$ws.Range("E3").Value = '  +1.54%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").Value = '414.97'
$ws.Range("E5").Value = '  +1.04%  '
$ws.Range("D6").Value = '129.35'
$ws.Range("E6").Value = '  -0.60%  '
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +0.38%  '
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("D9").Value = '0.756'
$ws.Range("E9").Value = '  +2.61%  '
$ws.Range("D10").Value = '0.156'
$ws.Range("E10").Value = '  +11.53%  '
$ws.Range("D11").Value = '42.49'
$ws.Range("E11").Value = '  -2.77%  '
$ws.Range("B12").Value = 'ShibaInu'
$ws.Range("C12").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D12").Value = '0.0000229'
$ws.Range("E12").Value = '  +1.17%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").Value = '9.68'
$ws.Range("E13").Value = '  +3.32%  '
$ws.Range("D14").Value = '4.029.73'
$ws.Range("E14").Value = '  +1.51%  '
$ws.Range("E15").Value = '  -0.91%  '
$ws.Range("D16").Value = '20.27'
$ws.Range("E16").Value = '  -4.67%  '
$ws.Range("D17").Value = '3.472.86'
$ws.Range("E17").Value = '  +1.58%  '
$ws.Range("D18").Value = '1.09'
$ws.Range("E18").Value = '  +1.01%  '
$ws.Range("D19").Value = '12.40'
$ws.Range("E19").Value = '  +0.06%  '
$ws.Range("D20").Value = '63.553.97'
$ws.Range("E20").Value = '  +2.64%  '
$ws.Range("D21").Value = '457.45'
$ws.Range("E21").Value = '  -10.93%  '
$ws.Range("D22").Value = '89.92'
$ws.Range("E22").Value = '  -2.81%  '
$ws.Range("E23").Value = '  -1.88%  '
$ws.Range("D24").Value = '13.19'
$ws.Range("E24").Value = '  -2.25%  '
$ws.Range("D25").Value = '10.18'
$ws.Range("E25").Value = '  +9.48%  '
$ws.Range("D26").Value = '3.31'
$ws.Range("E26").Value = '  -0.76%  '
$ws.Range("D27").Value = '33.50'
$ws.Range("E27").Value = '  -3.94%  '
$ws.Range("D28").Value = '4.76'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").Value = '12.49'
$ws.Range("E29").Value = '  +2.57%  '
$ws.Range("D30").Value = '7.52'
$ws.Range("E30").Value = '  -2.03%  '
$ws.Range("D31").Value = '2.67'
$ws.Range("E31").Value = '  -1.23%  '
$ws.Range("E32").Value = '  +1.18%  '
$ws.Range("D33").Value = '0.112'
$ws.Range("E33").Value = '  -2.10%  '
$ws.Range("D34").Value = '39.86'
$ws.Range("E34").Value = '  -5.19%  '
$ws.Range("D36").Value = '57.53'
$ws.Range("E36").Value = '  -3.34%  '
$ws.Range("D37").Value = '0.0485'
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("B38").Value = 'Stacks'
$ws.Range("C38").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D38").Value = '3.09'
$ws.Range("E38").Value = '  +4.83%  '
$ws.Range("B39").Value = 'FirstDigitalUSD'
$ws.Range("C39").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D39").Value = '0.999'
$ws.Range("E39").Value = '  +0.09%  '
$ws.Range("B40").Value = 'WEMIXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '2.81'
$ws.Range("E40").Value = '  +2.04%  '
$ws.Range("B41").Value = 'Stellar'
$ws.Range("C41").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D41").Value = '0.137'
$ws.Range("E41").Value = '  -1.66%  '
$ws.Range("D42").Value = '4.50'
$ws.Range("E42").Value = '  +4.02%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0647'
$ws.Range("E43").Value = '  +53.11%  '
$ws.Range("B44").Value = 'Monero'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D44").Value = '146.58'
$ws.Range("E44").Value = '  -0.94%  '
$ws.Range("B45").Value = 'LidoDAOToken'
$ws.Range("C45").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D45").Value = '3.32'
$ws.Range("E45").Value = '  -4.18%  '
$ws.Range("E46").Value = '  -1.39%  '
$ws.Range("D47").Value = '2.00'
$ws.Range("E47").Value = '  -5.90%  '
$ws.Range("D48").Value = '2.34'
$ws.Range("E48").Value = '  -1.15%  '
$ws.Range("D49").Value = '16.00'
$ws.Range("E49").Value = '  -4.48%  '
$ws.Range("D50").Value = '21.58'
$ws.Range("E50").Value = '  -6.72%  '
$ws.Range("D51").Value = '0.140'
$ws.Range("E51").Value = '  -5.29%  '

